# xls export geometry fix, added server power\health status
#
# Shortens a handful of inventory-report column headers so the exported
# sheet has room for the new power/health-status columns, and narrows the
# affected columns to match the shorter text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text shortenings -------------------------------------------------
$ws.Range("F1").Value = "Memory tot.size"   # was "System memory size"
$ws.Range("H1").Value = "Memory P/Ns"       # was "Memory module part number"
$ws.Range("M1").Value = "HDD slot pop."     # was "HDD slot population"
$ws.Range("N1").Value = "PSU P/Ns"          # was "PSU part number"

# --- Column width adjustments to match the shortened headers ----------------
# (values are the character-width inputs whose nearest representable column
# width in this engine matches the narrower geometry used by the export)
$ws.Columns.Item(6).ColumnWidth = 14.75    # column F: 18.7109375 -> ~15.71
$ws.Columns.Item(8).ColumnWidth = 10.75    # column H: 25.7109375 -> ~11.71
$ws.Columns.Item(13).ColumnWidth = 12.75   # column M: 19.7109375 -> ~13.71
$ws.Columns.Item(14).ColumnWidth = 7.75    # column N: 15.7109375 -> ~8.71
